# fix(publipostage): Refactor synthetic array /3
#
# The workbook uses colored-square emoji in column A ("statut") to denote a
# status, with a matching human readable color label in column B
# ("statut_label"). This change swaps the colored-square emoji for colored-
# book emoji, and renames the "noir" (black) label to "bleu" (blue) to match
# the new "blue book" emoji used for that status.
#
#   🟥 (rouge/red)    -> 📕
#   ⬛ (noir/black)    -> 📘  (label "noir" -> "bleu")
#   🟩 (vert/green)   -> 📗
#   🟧 (orange)       -> 📙

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: replace the status emoji -------------------------------

$redCells = @("A2","A3","A4","A8","A9","A10","A11","A14","A42","A65","A73","A98")
$blackCells = @("A5","A6","A7","A12","A13","A15","A16","A17","A19","A20","A21","A22","A23","A24","A25","A26","A27","A28","A29","A30","A31","A32","A34","A37","A38","A39","A41","A44","A46","A47","A48","A50","A52","A54","A55","A56","A57","A59","A60","A61","A62","A63","A64","A67","A68","A69","A70","A71","A74","A75","A76","A77","A78","A79","A80","A81","A82","A84","A85","A86","A87","A88","A89","A90","A91","A92","A93","A94","A95","A96","A97","A99","A100","A101")
$greenCells = @("A18","A33","A35","A40","A43","A49","A51","A58","A66","A72","A83")
$orangeCells = @("A36","A45","A53")

foreach ($c in $redCells)    { $ws.Range($c).Value = "📕" }
foreach ($c in $blackCells)  { $ws.Range($c).Value = "📘" }
foreach ($c in $greenCells)  { $ws.Range($c).Value = "📗" }
foreach ($c in $orangeCells) { $ws.Range($c).Value = "📙" }

# --- Column B: rename the "noir" label to "bleu" -----------------------

$noirCells = @("B5","B6","B7","B12","B13","B15","B16","B17","B19","B20","B21","B22","B23","B24","B25","B26","B27","B28","B29","B30","B31","B32","B34","B37","B38","B39","B41","B44","B46","B47","B48","B50","B52","B54","B55","B56","B57","B59","B60","B61","B62","B63","B64","B67","B68","B69","B70","B71","B74","B75","B76","B77","B78","B79","B80","B81","B82","B84","B85","B86","B87","B88","B89","B90","B91","B92","B93","B94","B95","B96","B97","B99","B100","B101")

foreach ($c in $noirCells) { $ws.Range($c).Value = "bleu" }
